$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_IE")

# Insert two new blank columns before column X (24), shifting existing
# X:AH content to Z:AJ. This matches the diff where two new columns
# (derived_variable / derivation_description) were added while all the
# former X..AH columns shifted right by two positions.
$ws.Range("X1:Y1").EntireColumn.Insert()

# Header labels for the two newly inserted columns.
$ws.Range("X1").Value = "derived_variable"
$ws.Range("Y1").Value = "derivation_description"

# Match the header cell formatting used by the other header cells in row 1.
$ws.Range("X1:Y1").Style = $ws.Range("W1").Style

# Update the selection / view to match the saved workbook state.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("W3").Select()
